$d = $word.ActiveDocument

# Replace "ooour first project" with "second commit."
$d.Content.Find.Execute("ooour first project", $true, $false, $false, $false, $false,
                         $true, 1, $false, "second commit.", 2)
